# Fix status name labels in the "statut_label" (column B) and
# "statut_name" (column C) columns of the publipostage workbook.
#
# Shared-string changes required:
#   "bleu"                                                  -> "noir"
#   "pas de résultat ni de publication"                     -> "pas de résultat postés ni publiés"
#   "résultat et / ou publication posté"                    -> "résultat postés ou publiés"
#   "résultat et / ou publication posté dans les 36 mois"   -> "résultat postés ou publiés dans les 36 mois"
#   "résultat et / ou publication posté dans les 12 mois"   -> "résultat postés ou publiés dans les 12 mois"
#
# Note: "résultat et / ou publication posté" is a textual prefix of the
# two "... dans les XX mois" variants, so replacements are done with a
# whole-cell match (xlWhole) to avoid corrupting the longer strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlWhole = 1          # Microsoft.Office.Interop.Excel.XlLookAt.xlWhole
$xlByRows = 1          # Microsoft.Office.Interop.Excel.XlSearchOrder.xlByRows
$xlNext = 1          # Microsoft.Office.Interop.Excel.XlSearchDirection.xlNext

$replacements = @(
    @{ Old = "bleu"; New = "noir" },
    @{ Old = "résultat et / ou publication posté dans les 36 mois"; New = "résultat postés ou publiés dans les 36 mois" },
    @{ Old = "résultat et / ou publication posté dans les 12 mois"; New = "résultat postés ou publiés dans les 12 mois" },
    @{ Old = "résultat et / ou publication posté"; New = "résultat postés ou publiés" },
    @{ Old = "pas de résultat ni de publication"; New = "pas de résultat postés ni publiés" }
)

foreach ($r in $replacements) {
    $ws.Cells.Replace(
        $r.Old,
        $r.New,
        $xlWhole,
        $xlByRows,
        $false,
        $false,
        $false,
        $false
    ) | Out-Null
}
